# Minor correction to slide 60 ("Maintaining Context During Parsing"):
# The sentence about the "exit" statement had a stray period before the
# comma ("...nested inside a loop., and code generation...").  Fix the
# typo so it reads "...nested inside a loop, and code generation...".
#
# The original text lives as a single run inside the 3rd paragraph of the
# body placeholder (Shapes.Item(2)) on slide 60.  We locate the run text
# dynamically (rather than relying on hard-coded character offsets) and
# replace just the broken fragment "a loop., " with the corrected
# "a loop, " using TextRange.Characters(start, length), which is the
# classic PowerPoint COM way to operate on a sub-range of characters.
# Editing only that inner fragment naturally splits the original single
# run into three runs -- matching the intended edit -- while leaving the
# surrounding text/formatting untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(60)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$oldFragment = "a loop., "
$newFragment = "a loop, "

$fullText = $tr.Text
$idx = $fullText.IndexOf($oldFragment)

if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $oldFragment.Length)
    $sub.Text = $newFragment
}
